$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.370.49"
$ws.Range("D3").Value = "1.604.99"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.831.32"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.610.51"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.28"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "26.370.44"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.46%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "1.489.23"
$ws.Range("E32").Value = "  +5.25%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.561"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.934"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").Value = "1.742.67"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.40%  "
